$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell reference -> new text value (kept as plain text, matching
# the workbook's existing inline-string storage for these columns).
$updates = @{
    'D2' = '303.46'
    'E2' = '5.00%'
    'D3' = '32.12'
    'E3' = '9.81%'
    'D4' = '5.258'
    'E4' = '0.05%'
    'D5' = '0.07514'
    'E5' = '6.48%'
    'D6' = '7.870'
    'E6' = '5.50%'
    'D7' = '3.817'
    'E7' = '7.00%'
    'D8' = '1.514'
    'E8' = '8.24%'
    'D9' = '0.9213'
    'E9' = '1.71%'
    'D10' = '0.1704'
    'E10' = '5.75%'
    'D11' = '0.07903'
    'E11' = '3.88%'
    'D12' = '0.08030'
    'E12' = '3.03%'
    'D13' = '0.03046'
    'E13' = '4.66%'
    'D14' = '0.09909'
    'E14' = '9.81%'
    'D15' = '0.001489'
    'E15' = '-6.32%'
    'D16' = '0.04607'
    'E16' = '1.98%'
    'D17' = '0.006343'
    'E17' = '2.70%'
    'D18' = '3.462'
    'E18' = '-0.70%'
    'D19' = '2.233'
    'E19' = '0.00%'
    'D20' = '0.3304'
    'E20' = '2.19%'
    'E21' = '-0.53%'
    'D22' = '4.505'
    'E22' = '11.75%'
    'E23' = '1.11%'
    'D24' = '0.001216'
    'E24' = '0.47%'
    'D25' = '0.004448'
    'E25' = '6.73%'
    'D26' = '0.0001397'
    'E26' = '19.64%'
    'D27' = '0.0001773'
    'E27' = '6.37%'
    'E39' = '2,549.28%'
    'D40' = '0.04487'
    'E40' = '2.27%'
    'D41' = '0.006959'
    'E41' = '-0.76%'
    'D42' = '0.1352'
    'E42' = '7.99%'
    'D43' = '0.002196'
    'E43' = '6.26%'
    'D44' = '0.01285'
    'E44' = '7.80%'
    'D45' = '0.00006163'
    'E45' = '5.61%'
    'D47' = '0.01497'
    'E47' = '15.37%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Force text format so numeric-looking / percent-looking strings are not
    # reinterpreted by Excel as numbers (the source cells are plain text).
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}
